$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A, D, J move from "dimension" to "measure" for
# grupo-de-tipo-de-jornada, mes-nombre and sexo respectively.

# Row 2: update the iaest-dimension:... -> iaest-measure:... labels
$ws.Range("A2").Value = "iaest-measure:grupo-de-tipo-de-jornada"
$ws.Range("D2").Value = "iaest-measure:mes-nombre"
$ws.Range("J2").Value = "iaest-measure:sexo"

# Row 3: "dim" -> "medida" for the same columns
$ws.Range("A3").Value = "medida"
$ws.Range("D3").Value = "medida"
$ws.Range("J3").Value = "medida"

# Row 4: "skos:Concept" -> "xsd:int" for the same columns
$ws.Range("A4").Value = "xsd:int"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"

# Row 5: mapping file references no longer apply to these
# now-measure columns, so the cells are removed entirely.
$ws.Range("A5").Clear()
$ws.Range("D5").Clear()
$ws.Range("J5").Clear()
